$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "2011-11-21"
$ws.Range("J2").NumberFormat = "General"
Write-Output "done"
